$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns AD, AE, AF ("Wins", "Losses", "Ties")
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting/style of the existing header row (copy from A1)
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in Wins/Losses/Ties values for each data row (rows 2-45)
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 78  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 84  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
